# Added width of columns
#
# Sheet1: turn row 1 into a header row ("COL 1" / "COL 2"), keep the
# existing "Sample" cell, and append a couple of new data cells/rows.
# Sheet2: a brand-new sheet with the same header row plus one data cell.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 -----------------------------------------------------------
# A1 keeps its original (bordered / centered) style, now holding text.
$ws1.Range("A1").Value = "COL 1"

# B1 is a new bold header cell (no border).
$ws1.Range("B1").Value = "COL 2"
$ws1.Range("B1").Font.Bold = $true

# Extra data below the header row.
$ws1.Range("B2").Value = "Very large data value present here"
$ws1.Range("A3").Value = "New"

# --- Sheet2 (new) -------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "COL 1"
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("B1").Value = "COL 2"
$ws2.Range("B1").Font.Bold = $true

$ws2.Range("B2").Value = "data"

# --- View state: restore each sheet's selection, Sheet2 ends up active -
$ws1.Range("D8").Select() | Out-Null
$ws2.Range("C4").Select() | Out-Null
